$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.298330068588257
$ws.Range("B1").Value = 2.47739577293396
$ws.Range("C1").Value = 3.175394058227539
$ws.Range("D1").Value = 1.586362957954407
$ws.Range("E1").Value = 1.145838856697083
